# Generate Report for Handoff
#
# Updates the "9924594f-349d-4e42-bdc4-7bd1aa074a90" row's handoff/generate
# timestamps across all three sheets now that the handoff xliff/report was
# (re)generated:
#   - Overview!G7             "Latest HO Xliff Generate Date" -> 2016-08-28 00:39:30
#   - zh-cn!H7                "Latest Handoff Datetime"       -> 2016-08-28 00:39:26
#   - de-de!H7                "Latest Handoff Datetime"       -> 2016-08-28 00:39:30

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G7").Value = "2016-08-28 00:39:30"

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H7").Value = "2016-08-28 00:39:26"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H7").Value = "2016-08-28 00:39:30"
